$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each target cell stores a plain text value (prices, percentages and the
# hour marker are text in this sheet, not numbers), so force text storage
# via NumberFormat "@" immediately before writing each value - this keeps
# numeric-looking strings (e.g. "308.28", "-1.07%", "17") from being
# auto-converted into real numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.07%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "17"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.38"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.10%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "17"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.135"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.05%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "17"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07723"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.19%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "17"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.56%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "17"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.294"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.33%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "17"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.846"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.61%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "17"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.965"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-4.22%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "17"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9210"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.77%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "17"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1130"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-7.96%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "17"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1853"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.78%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "17"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08751"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.00%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "17"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.30%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "17"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09533"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.04%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "17"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001378"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.18%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "17"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006011"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "5.59%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "17"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.374"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.25%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "17"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3447"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.39%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "17"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.324"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "19.28%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "17"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1299"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.08%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "17"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-10.62%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "17"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04340"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.71%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "17"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-3.13%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "17"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004263"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.94%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "17"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001333"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "9.22%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "17"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002904"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "17"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "17"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "17"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "17"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "17"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "17"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "17"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "17"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "17"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "17"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "17"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "17"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.03%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "17"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04906"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-5.06%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "17"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007580"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.02%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "17"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1350"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.80%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "17"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008541"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-6.58%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "17"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.57%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "17"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008546"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.87%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "17"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006428"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.16%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "17"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.19%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "17"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "17.73%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "17"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001445"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "20.34%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "17"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.19%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "17"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.19%"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "17"
